$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  8"
$ws.Range("C9").Value = "Report Covering the Week  2/17/2025  Through  2/23/2025"

# --- Data cell updates (rows 14-33) ---
# Row 14
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "***.*"
$ws.Range("F14").Copy()
$ws.Range("H14").PasteSpecial(-4122)

# Row 15
$ws.Range("C16").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 2
$ws.Range("C16").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 2
$ws.Range("E16").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 4
$ws.Range("C16").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = 2
$ws.Range("E16").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 5
$ws.Range("C16").Copy()
$ws.Range("J15").PasteSpecial(-4122)
$ws.Range("J15").Value = 2
$ws.Range("E16").Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("K15").Value = 150
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 66.666666666666
$ws.Range("N15").Value = -16.666666666666

# Row 16
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 19
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = -32.142857142857
$ws.Range("I16").Value = 27
$ws.Range("J16").Value = 53
$ws.Range("K16").Value = -49.056603773584
$ws.Range("L16").Value = -59.090909090909
$ws.Range("M16").Value = -32.5
$ws.Range("N16").Value = -85.483870967741

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = 19.230769230769
$ws.Range("I17").Value = 62
$ws.Range("J17").Value = 58
$ws.Range("K17").Value = 6.896551724137
$ws.Range("L17").Value = 1.639344262295
$ws.Range("M17").Value = 244.444444444444
$ws.Range("N17").Value = 44.186046511627

# Row 18
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 37
$ws.Range("H18").Value = 15.625
$ws.Range("I18").Value = 80
$ws.Range("J18").Value = 66
$ws.Range("K18").Value = 21.212121212121
$ws.Range("L18").Value = -20
$ws.Range("M18").Value = 3.896103896103
$ws.Range("N18").Value = -78.891820580474

# Row 19
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 26
$ws.Range("E19").Value = -19.230769230769
$ws.Range("F19").Value = 58
$ws.Range("G19").Value = 108
$ws.Range("H19").Value = -46.296296296296
$ws.Range("I19").Value = 110
$ws.Range("J19").Value = 207
$ws.Range("K19").Value = -46.859903381642
$ws.Range("L19").Value = -50
$ws.Range("M19").Value = 41.025641025641
$ws.Range("N19").Value = -41.489361702127

# Row 20
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 32
$ws.Range("G20").Value = 34
$ws.Range("H20").Value = -5.882352941176
$ws.Range("I20").Value = 61
$ws.Range("J20").Value = 71
$ws.Range("K20").Value = -14.084507042253
$ws.Range("L20").Value = -7.575757575757
$ws.Range("M20").Value = 60.526315789473
$ws.Range("N20").Value = -90.113452188006

# Row 21
$ws.Range("C21").Value = 50
$ws.Range("D21").Value = 55
$ws.Range("E21").Value = -9.090909090909
$ws.Range("F21").Value = 181
$ws.Range("G21").Value = 230
$ws.Range("H21").Value = -21.304347826087
$ws.Range("I21").Value = 345
$ws.Range("J21").Value = 458
$ws.Range("K21").Value = -24.672489082969
$ws.Range("L21").Value = -33.397683397683
$ws.Range("M21").Value = 35.294117647058
$ws.Range("N21").Value = -75.704225352112

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("F14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 5
$ws.Range("K22").Value = 66.666666666666
$ws.Range("L22").Value = 66.666666666666

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("F14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("F23").Value = 2
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 2
$ws.Range("K23").Value = -33.333333333333
$ws.Range("L23").Value = -33.333333333333
$ws.Range("M23").Value = 100

# Row 24
$ws.Range("C24").Value = 61
$ws.Range("D24").Value = 56
$ws.Range("E24").Value = 8.928571428571
$ws.Range("F24").Value = 206
$ws.Range("G24").Value = 231
$ws.Range("H24").Value = -10.82251082251
$ws.Range("I24").Value = 366
$ws.Range("J24").Value = 443
$ws.Range("K24").Value = -17.381489841986
$ws.Range("L24").Value = -12.649164677804
$ws.Range("M24").Value = 64.864864864864

# Row 25
$ws.Range("C25").Value = 38
$ws.Range("D25").Value = 37
$ws.Range("E25").Value = 2.702702702702
$ws.Range("F25").Value = 124
$ws.Range("G25").Value = 139
$ws.Range("H25").Value = -10.791366906474
$ws.Range("I25").Value = 215
$ws.Range("J25").Value = 280
$ws.Range("K25").Value = -23.214285714285
$ws.Range("L25").Value = -1.37614678899

# Row 26
$ws.Range("C26").Value = 23
$ws.Range("D26").Value = 20
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = 65
$ws.Range("G26").Value = 51
$ws.Range("H26").Value = 27.450980392156
$ws.Range("I26").Value = 128
$ws.Range("J26").Value = 104
$ws.Range("K26").Value = 23.076923076923
$ws.Range("L26").Value = 8.474576271186
$ws.Range("M26").Value = 50.588235294117

# Row 27
$ws.Range("C16").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 2
$ws.Range("C16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 2
$ws.Range("E16").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 5
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = 25
$ws.Range("L27").Value = -37.5

# Row 28
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 11
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = 22.222222222222
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 12.5
$ws.Range("L28").Value = 80

# Row 29
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("C29").PasteSpecial(-4122)

# Row 30
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("C30").PasteSpecial(-4122)

# Row 31
$ws.Range("E16").Copy()
$ws.Range("L31").PasteSpecial(-4122)
$ws.Range("L31").Value = -100

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "***.*"
$ws.Range("F14").Copy()
$ws.Range("E33").PasteSpecial(-4122)

$excel.CutCopyMode = 0
Write-Output "done"